# Updated the location in excel
$wb = $excel.ActiveWorkbook

$wsLogin = $wb.Worksheets.Item("UserLogin")
$wsFiles = $wb.Worksheets.Item("FilesToUpload")
$wsUpload = $wb.Worksheets.Item("FIles to upload")

$networkPath = "\\ehs-clu-bos-001.ehs.govt.state.ma.us\File Services\QA\Projects\VG\FTS Automation\Diane - test files\Original"

# ----- UserLogin sheet -----
# Header row
$wsLogin.Range("A1").Value = "Env"
$wsLogin.Range("B1").Value = "UserName"
$wsLogin.Range("C1").Value = "Password"
$wsLogin.Range("D1").Value = "User"
$wsLogin.Range("E1").Value = "Location"
$wsLogin.Range("F1").Value = "Files"
$wsLogin.Range("A1:F1").Font.Bold = $true

# Data rows (3 environments: bphc, sinc1, hpoint)
$wsLogin.Range("A2").Value = "QA"
$wsLogin.Range("B2").Value = "bphc"
$wsLogin.Range("C2").Value = "Ehsquincyqa21!!"
$wsLogin.Range("D2").Value = "bphc"
$wsLogin.Range("E2").Value = ""

$wsLogin.Range("A3").Value = "QA"
$wsLogin.Range("B3").Value = "sinc1"
$wsLogin.Range("C3").Value = "EhsBostonQA22!!"
$wsLogin.Range("D3").Value = "sinc1"
$wsLogin.Range("E3").Value = ""

$wsLogin.Range("A4").Value = "QA"
$wsLogin.Range("B4").Value = "hpoint"
$wsLogin.Range("C4").Value = "Ehsquincyqa20!!"
$wsLogin.Range("D4").Value = "hpoint"
$wsLogin.Range("E4").Value = "High Point Treatment Ctr"

# Remove the old trailing rows (table shrank from 7 to 4 rows)
$wsLogin.Rows.Item(7).Delete()
$wsLogin.Rows.Item(6).Delete()
$wsLogin.Rows.Item(5).Delete()

# Files column becomes a hyperlink to the shared network location
$wsLogin.Hyperlinks.Add($wsLogin.Range("F2"), $networkPath)
$wsLogin.Hyperlinks.Add($wsLogin.Range("F3"), $networkPath)
$wsLogin.Hyperlinks.Add($wsLogin.Range("F4"), $networkPath)

# Column widths
$wsLogin.Columns.Item(5).ColumnWidth = 21.6328125
$wsLogin.Columns.Item(6).ColumnWidth = 95.36328125

# ----- FilesToUpload sheet -----
$wsFiles.Range("A1").Value = "bphc_FileName"
$wsFiles.Range("B1").Value = "bphc_FileType"
$wsFiles.Range("C1").Value = "sinc1_FileName"
$wsFiles.Range("D1").Value = "sinc1_FileType"
$wsFiles.Range("E1").Value = "hpoint_FileName"
$wsFiles.Range("F1").Value = "hpoint_FileType"

$wsFiles.Range("A2").Value = "LARGEBPC.HIPAA937P"
$wsFiles.Range("B2").Value = "PROF"
$wsFiles.Range("C2").Value = "LargeFileSpectrum.HIPAA837P"
$wsFiles.Range("D2").Value = "PROF"
$wsFiles.Range("E2").Value = "LargeHighPointTreatCtr.HIPAA837I"
$wsFiles.Range("F2").Value = "PROF"
$wsFiles.Range("A2:F2").Font.Bold = $false

# Remove old trailing rows (table shrank from 7 to 2 rows)
$wsFiles.Rows.Item(7).Delete()
$wsFiles.Rows.Item(6).Delete()
$wsFiles.Rows.Item(5).Delete()
$wsFiles.Rows.Item(4).Delete()
$wsFiles.Rows.Item(3).Delete()

$wsFiles.Columns.Item(1).ColumnWidth = 37.6328125
$wsFiles.Columns.Item(2).ColumnWidth = 16.08984375
$wsFiles.Columns.Item(3).ColumnWidth = 34.36328125
$wsFiles.Columns.Item(4).ColumnWidth = 17.1796875
$wsFiles.Columns.Item(5).ColumnWidth = 34.1796875

# ----- View / selection bookkeeping -----
$wsUpload.Range("B24:C29").Select()
$wsFiles.Range("F6").Select()
$wsLogin.Range("E11").Select()
